$d = $word.ActiveDocument

$replacements = @(
    @("500×6=", "374×3="),
    @("393×9=", "104×3="),
    @("626×8=", "198×3="),
    @("368×4=", "251×5="),
    @("766×2=", "463×5="),
    @("820×7=", "550×5="),
    @("139×2=", "387×5="),
    @("586×6=", "554×6="),
    @("876×9=", "604×7="),
    @("635×5=", "783×7="),
    @("873×7=", "974×5="),
    @("455×4=", "110×2="),
    @("504×5=", "311×4="),
    @("150×8=", "115×4="),
    @("309×7=", "542×6="),
    @("933×9=", "964×9="),
    @("721×5=", "414×3="),
    @("761×8=", "230×7="),
    @("611×8=", "826×5="),
    @("721×7=", "120×5="),
    @("223×5=", "419×8="),
    @("928×7=", "132×3="),
    @("365×2=", "554×4="),
    @("448×5=", "837×4="),
    @("354×5=", "555×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
